$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Units" column (N) values on the data rows (2-6) no longer apply
# for these ordinal/nominal scale rows - clear them out.
$ws.Range("N2:N6").ClearContents()

# Reflect the last-touched range as the active selection, as Excel does
# after an edit.
$ws.Range("N2:N6").Select()
